$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.93000030517578
$ws.Range("F2").Value = 0.21
$ws.Range("I2").Value = 0.239999994635582
$ws.Range("L2").Value = "19.93±0.22"
$ws.Range("O2").Value = 87.67
$ws.Range("C3").Value = 25.29000091552734
$ws.Range("F3").Value = 0.3
$ws.Range("L3").Value = "25.29±0.22"
$ws.Range("O3").Value = 84.61
$ws.Range("C4").Value = 30.09000015258789
$ws.Range("F4").Value = 0.11
$ws.Range("L4").Value = "30.09±0.1"
$ws.Range("O4").Value = 106.83
$ws.Range("C5").Value = 31.70999908447266
$ws.Range("F5").Value = 0.29
$ws.Range("L5").Value = "31.71±0.12"
$ws.Range("O5").Value = 91.41
$ws.Range("C6").Value = 34.88000106811523
$ws.Range("F6").Value = 0.13
$ws.Range("L6").Value = "34.88±0.1"
$ws.Range("O6").Value = 107.65
$ws.Range("C7").Value = 37.16999816894531
$ws.Range("F7").Value = 0.17
$ws.Range("I7").Value = 0.1000000014901161
$ws.Range("L7").Value = "37.17±0.09"
$ws.Range("O7").Value = 104.74
$ws.Range("C8").Value = 39.93999862670898
$ws.Range("L8").Value = "39.94±0.08"
$ws.Range("O8").Value = 115.33
$ws.Range("C9").Value = 43.02000045776367
$ws.Range("L9").Value = "43.02±0.1"
$ws.Range("O9").Value = 115.02
$ws.Range("C10").Value = 45.04999923706055
$ws.Range("L10").Value = "45.05±0.19"
$ws.Range("O10").Value = 105.55
$ws.Range("C11").Value = 47.02999877929688
$ws.Range("L11").Value = "47.03±0.16"
$ws.Range("O11").Value = 110.41
$ws.Range("C12").Value = 53.09999847412109
$ws.Range("F12").Value = 0.18
$ws.Range("L12").Value = "53.1±0.17"
$ws.Range("O12").Value = 107.92
$ws.Range("C13").Value = 55.06999969482422
$ws.Range("F13").Value = 0.15
$ws.Range("L13").Value = "55.07±0.15"
$ws.Range("O13").Value = 113.7
$ws.Range("C14").Value = 56.79000091552734
$ws.Range("F14").Value = 0.25
$ws.Range("L14").Value = "56.79±0.21"
$ws.Range("O14").Value = 104.44
$ws.Range("C15").Value = 59.84000015258789
$ws.Range("F15").Value = 0.22
$ws.Range("L15").Value = "59.84±0.18"
$ws.Range("O15").Value = 108.81
$ws.Range("C16").Value = 63.18000030517578
$ws.Range("F16").Value = 0.24
$ws.Range("L16").Value = "63.18±0.18"
$ws.Range("O16").Value = 106.91
$ws.Range("C17").Value = 65.20999908447266
$ws.Range("F17").Value = 0.25
$ws.Range("L17").Value = "65.21±0.22"
$ws.Range("O17").Value = 104.07
$ws.Range("C18").Value = 67.05000305175781
$ws.Range("F18").Value = 0.2
$ws.Range("L18").Value = "67.05±0.2"
$ws.Range("O18").Value = 110.5
$ws.Range("C19").Value = 69.87000274658203
$ws.Range("F19").Value = 0.24
$ws.Range("L19").Value = "69.87±0.22"
$ws.Range("O19").Value = 110.35
$ws.Range("C20").Value = 72.95999908447266
$ws.Range("F20").Value = 0.23
$ws.Range("L20").Value = "72.96±0.22"
$ws.Range("O20").Value = 110.47
